$d = $word.ActiveDocument

# Locate the paragraph that ends with "...No need to change anything." so we
# can insert the new bullet right after it (and before the trailing, empty
# list paragraph).
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*method and select Create Unit Tests*") {
        $target = $p
        break
    }
}

# Create a new paragraph after the target one; it inherits the same
# ListParagraph style / numPr (numId 1, ilvl 0) from $target automatically.
$target.Range.InsertParagraphAfter() | Out-Null
$newP = $target.Next()
$newRange = $newP.Range

# Build the new bullet's content via a WordprocessingML package fragment so
# that the "NumberFunTests" word is wrapped in proofErr spellStart/spellEnd
# markers (matching how Word marks an unrecognized word), the same way the
# rest of the document does for its own flagged terms, while still keeping
# the paragraph's List Paragraph style/numbering.
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Add Arrange, Act, and Assert to the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>NumberFunTests</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>'' BiggestTestCase1.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$newRange.InsertXML($xml)
